$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Fix typo in FC1 description (B3): "tactil" -> "tactile"
$ws.Range("B3").Value = "Traiter les informations provenant de l'écran tactile et des boutons"

# Add description for FC2 (B4): new function description
$ws.Range("B4").Value = "S'intégrer au boitier en n'altérant pas le design de l'objet"

# Update the active selection to B4
$ws.Range("B4").Select()
